$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Valeurs de cellules (pieuvre des fonctions : ajout du phasage FC3, mise a
# jour des criteres/niveaux d'exigence FC1/FC4/FC5 et de la note de bas de
# page des normes). L'ordre d'ecriture ci-dessous respecte l'ordre
# d'apparition des nouvelles chaines partagees (sharedStrings) telles que
# produites par le classeur cible, pour que la table compactee corresponde.
# ---------------------------------------------------------------------------

$ws.Range("D3").Value = "Moins de 100 ms entre l'appui et la réponse"
$ws.Range("E3").Value = "F1`nLatence : ±20ms"
$ws.Range("B5").Value = "Etablir une connectivité et traiter les informations provenant de l'application compagnon"
$ws.Range("D5").Value = "Moins de 100 ms entre l'appui et la réponse`nPortée de 40m en indoor"
$ws.Range("C5").Value = "Faible latence et fluidité des contrôles`nPortée importante`nMinimalisme des commandes possibles"
$ws.Range("C3").Value = "Faible latence et fluidité des contrôles`nMinimalisme des commandes possibles"
$ws.Range("B6").Value = "Respect des différentes normes en vigueur*"
$ws.Range("A9").Value = "*"
$ws.Range("B9").Value = "Normes CE principalement : `nCompatibilité électromagnétique (CEM) - 2014/30/UE`nÉquipements terminaux de télécommunication - 1999/5/CE "
$ws.Range("D6").Value = "Toutes les normes doivent être respectées"
$ws.Range("C6").Value = "Respect des normes"
$ws.Range("E6").Value = "F0"
$ws.Range("E7").Value = "F0"
$ws.Range("B7").Value = "Traiter les données provenant du tuner"
$ws.Range("C7").Value = "Traitement rapide et sans perte"
$ws.Range("D7").Value = "Pas de pertes de données"
$ws.Range("D2").Value = "2x3W, son 'chaud' et non saturé"
$ws.Range("E5").Value = "F1`nLatence : ±20ms`nPortée : ±15m"

# Cellules dont le texte ne change pas mais dont la mise en forme change
$ws.Range("B3").Value = "Traiter les informations provenant de l'écran tactile et des boutons"
$ws.Range("B4").Value = "S'intégrer au boitier en n'altérant pas le design de l'objet"

# ---------------------------------------------------------------------------
# Mise en forme : habillage du texte (wrap) -> style 1
# ---------------------------------------------------------------------------
$ws.Range("E2").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("E3").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("D5").WrapText = $true
$ws.Range("E5").WrapText = $true
$ws.Range("E6").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("C7").WrapText = $true
$ws.Range("B9").WrapText = $true

# ---------------------------------------------------------------------------
# Mise en forme : alignement vertical centre -> style 2
# ---------------------------------------------------------------------------
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Hauteur des lignes
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# ---------------------------------------------------------------------------
# Largeur des colonnes
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 80.41666666666683
$ws.Columns.Item(3).ColumnWidth = 35.91666666666654
$ws.Columns.Item(4).ColumnWidth = 39.08333333333342

# ---------------------------------------------------------------------------
# Selection active
# ---------------------------------------------------------------------------
$ws.Range("B9").Select()
